# Update Name of Algo
# Apply revised KNN-imputed values in column A for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value  = -20.257
$ws.Range("A6").Value  = -22.291
$ws.Range("A7").Value  = -19.851
$ws.Range("A16").Value = -21.93
$ws.Range("A20").Value = -20.137
$ws.Range("A28").Value = -21.893
$ws.Range("A29").Value = -21.344
$ws.Range("A32").Value = -21.648
$ws.Range("A40").Value = -19.896
$ws.Range("A46").Value = -21.87
$ws.Range("A51").Value = -21.98
$ws.Range("A52").Value = -22.084
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.395
$ws.Range("A62").Value = -22.085
$ws.Range("A66").Value = -21.616
$ws.Range("A73").Value = -20.213
$ws.Range("A74").Value = -21.112
$ws.Range("A92").Value = -21.50600000000001
$ws.Range("A100").Value = -22.217
